$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 16,20
$arr[0,0] = "ECs"
$arr[0,1] = "Reln"
$arr[0,2] = "Itgb1"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 0.06455366666666666
$arr[0,7] = 0.193661
$arr[0,8] = 0.01357839286814829
$arr[0,9] = 0.01357839286814829
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 168.1098273333333
$arr[0,13] = 504.329482
$arr[0,14] = 0.2984182258032519
$arr[0,15] = 0.298418225803252
$arr[0,16] = 10.85210575706689
$arr[0,17] = 97.66895181360199
$arr[0,18] = 0.004052039908972341
$arr[0,19] = 0.004052039908972343
$arr[1,0] = "ECs"
$arr[1,1] = "Reln"
$arr[1,2] = "Itgb1"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 0.06455366666666666
$arr[1,7] = 0.193661
$arr[1,8] = 0.01357839286814829
$arr[1,9] = 0.01357839286814829
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 163.0062356666667
$arr[1,13] = 489.018707
$arr[1,14] = 0.2893586437755394
$arr[1,15] = 0.2893586437755394
$arr[1,16] = 10.52265020181411
$arr[1,17] = 94.703851816327
$arr[1,18] = 0.003929025344978846
$arr[1,19] = 0.003929025344978846
$arr[2,0] = "ECs"
$arr[2,1] = "Reln"
$arr[2,2] = "Itgb1"
$arr[2,3] = "MuSCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 0.06455366666666666
$arr[2,7] = 0.193661
$arr[2,8] = 0.01357839286814829
$arr[2,9] = 0.01357839286814829
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 165.99353
$arr[2,13] = 497.98059
$arr[2,14] = 0.294661504941043
$arr[2,15] = 0.294661504941043
$arr[2,16] = 10.71549100444333
$arr[2,17] = 96.43941903999
$arr[2,18] = 0.0040010296772093
$arr[2,19] = 0.0040010296772093
$arr[3,0] = "ECs"
$arr[3,1] = "Reln"
$arr[3,2] = "Itgb1"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 0.06455366666666666
$arr[3,7] = 0.193661
$arr[3,8] = 0.01357839286814829
$arr[3,9] = 0.01357839286814829
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 66.22673433333334
$arr[3,13] = 198.680203
$arr[3,14] = 0.1175616254801657
$arr[3,15] = 0.1175616254801657
$arr[3,16] = 4.275178532575889
$arr[3,17] = 38.476606793183
$arr[3,18] = 0.001596297936987801
$arr[3,19] = 0.001596297936987802
$arr[4,0] = "FAPs"
$arr[4,1] = "Reln"
$arr[4,2] = "Itgb1"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.06084700000000001
$arr[4,7] = 0.182541
$arr[4,8] = 0.01279872257472933
$arr[4,9] = 0.01279872257472933
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 168.1098273333333
$arr[4,13] = 504.329482
$arr[4,14] = 0.2984182258032519
$arr[4,15] = 0.298418225803252
$arr[4,16] = 10.22897866375133
$arr[4,17] = 92.060807973762
$arr[4,18] = 0.003819372083298756
$arr[4,19] = 0.003819372083298756
$arr[5,0] = "FAPs"
$arr[5,1] = "Reln"
$arr[5,2] = "Itgb1"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 0.06084700000000001
$arr[5,7] = 0.182541
$arr[5,8] = 0.01279872257472933
$arr[5,9] = 0.01279872257472933
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 163.0062356666667
$arr[5,13] = 489.018707
$arr[5,14] = 0.2893586437755394
$arr[5,15] = 0.2893586437755394
$arr[5,16] = 9.918440421609668
$arr[5,17] = 89.265963794487
$arr[5,18] = 0.00370342100628306
$arr[5,19] = 0.003703421006283059
$arr[6,0] = "FAPs"
$arr[6,1] = "Reln"
$arr[6,2] = "Itgb1"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 0.06084700000000001
$arr[6,7] = 0.182541
$arr[6,8] = 0.01279872257472933
$arr[6,9] = 0.01279872257472933
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 165.99353
$arr[6,13] = 497.98059
$arr[6,14] = 0.294661504941043
$arr[6,15] = 0.294661504941043
$arr[6,16] = 10.10020831991
$arr[6,17] = 90.90187487919
$arr[6,18] = 0.003771290855192645
$arr[6,19] = 0.003771290855192645
$arr[7,0] = "FAPs"
$arr[7,1] = "Reln"
$arr[7,2] = "Itgb1"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 0.06084700000000001
$arr[7,7] = 0.182541
$arr[7,8] = 0.01279872257472933
$arr[7,9] = 0.01279872257472933
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 66.22673433333334
$arr[7,13] = 198.680203
$arr[7,14] = 0.1175616254801657
$arr[7,15] = 0.1175616254801657
$arr[7,16] = 4.029698103980334
$arr[7,17] = 36.267282935823
$arr[7,18] = 0.001504638629954871
$arr[7,19] = 0.001504638629954871
$arr[8,0] = "MuSCs"
$arr[8,1] = "Reln"
$arr[8,2] = "Itgb1"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 4.62452
$arr[8,7] = 13.87356
$arr[8,8] = 0.9727340463997778
$arr[8,9] = 0.9727340463997778
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 168.1098273333333
$arr[8,13] = 504.329482
$arr[8,14] = 0.2984182258032519
$arr[8,15] = 0.298418225803252
$arr[8,16] = 777.4272586995465
$arr[8,17] = 6996.84532829592
$arr[8,18] = 0.2902815683050398
$arr[8,19] = 0.2902815683050399
$arr[9,0] = "MuSCs"
$arr[9,1] = "Reln"
$arr[9,2] = "Itgb1"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 4.62452
$arr[9,7] = 13.87356
$arr[9,8] = 0.9727340463997778
$arr[9,9] = 0.9727340463997778
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 163.0062356666667
$arr[9,13] = 489.018707
$arr[9,14] = 0.2893586437755394
$arr[9,15] = 0.2893586437755394
$arr[9,16] = 753.8255969652132
$arr[9,17] = 6784.43037268692
$arr[9,18] = 0.2814690044205324
$arr[9,19] = 0.2814690044205324
$arr[10,0] = "MuSCs"
$arr[10,1] = "Reln"
$arr[10,2] = "Itgb1"
$arr[10,3] = "MuSCs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 4.62452
$arr[10,7] = 13.87356
$arr[10,8] = 0.9727340463997778
$arr[10,9] = 0.9727340463997778
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 165.99353
$arr[10,13] = 497.98059
$arr[10,14] = 0.294661504941043
$arr[10,15] = 0.294661504941043
$arr[10,16] = 767.6403993555999
$arr[10,17] = 6908.7635942004
$arr[10,18] = 0.2866272780195488
$arr[10,19] = 0.2866272780195488
$arr[11,0] = "MuSCs"
$arr[11,1] = "Reln"
$arr[11,2] = "Itgb1"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 4.62452
$arr[11,7] = 13.87356
$arr[11,8] = 0.9727340463997778
$arr[11,9] = 0.9727340463997778
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 66.22673433333334
$arr[11,13] = 198.680203
$arr[11,14] = 0.1175616254801657
$arr[11,15] = 0.1175616254801657
$arr[11,16] = 306.2668574591867
$arr[11,17] = 2756.40171713268
$arr[11,18] = 0.1143561956546568
$arr[11,19] = 0.1143561956546568
$arr[12,0] = "Resolving-Mac"
$arr[12,1] = "Reln"
$arr[12,2] = "Itgb1"
$arr[12,3] = "ECs"
$arr[12,4] = 1
$arr[12,5] = 0.3333333333333333
$arr[12,6] = 0.004225666666666667
$arr[12,7] = 0.012677
$arr[12,8] = 0.0008888381573446169
$arr[12,9] = 0.000888838157344617
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 168.1098273333333
$arr[12,13] = 504.329482
$arr[12,14] = 0.2984182258032519
$arr[12,15] = 0.298418225803252
$arr[12,16] = 0.7103760937015555
$arr[12,17] = 6.393384843314
$arr[12,18] = 0.0002652455059410123
$arr[12,19] = 0.0002652455059410123
$arr[13,0] = "Resolving-Mac"
$arr[13,1] = "Reln"
$arr[13,2] = "Itgb1"
$arr[13,3] = "FAPs"
$arr[13,4] = 1
$arr[13,5] = 0.3333333333333333
$arr[13,6] = 0.004225666666666667
$arr[13,7] = 0.012677
$arr[13,8] = 0.0008888381573446169
$arr[13,9] = 0.000888838157344617
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 163.0062356666667
$arr[13,13] = 489.018707
$arr[13,14] = 0.2893586437755394
$arr[13,15] = 0.2893586437755394
$arr[13,16] = 0.6888100165154445
$arr[13,17] = 6.199290148639
$arr[13,18] = 0.0002571930037451879
$arr[13,19] = 0.0002571930037451879
$arr[14,0] = "Resolving-Mac"
$arr[14,1] = "Reln"
$arr[14,2] = "Itgb1"
$arr[14,3] = "MuSCs"
$arr[14,4] = 1
$arr[14,5] = 0.3333333333333333
$arr[14,6] = 0.004225666666666667
$arr[14,7] = 0.012677
$arr[14,8] = 0.0008888381573446169
$arr[14,9] = 0.000888838157344617
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 165.99353
$arr[14,13] = 497.98059
$arr[14,14] = 0.294661504941043
$arr[14,15] = 0.294661504941043
$arr[14,16] = 0.7014333266033334
$arr[14,17] = 6.31289993943
$arr[14,18] = 0.0002619063890921884
$arr[14,19] = 0.0002619063890921884
$arr[15,0] = "Resolving-Mac"
$arr[15,1] = "Reln"
$arr[15,2] = "Itgb1"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 1
$arr[15,5] = 0.3333333333333333
$arr[15,6] = 0.004225666666666667
$arr[15,7] = 0.012677
$arr[15,8] = 0.0008888381573446169
$arr[15,9] = 0.000888838157344617
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 66.22673433333334
$arr[15,13] = 198.680203
$arr[15,14] = 0.1175616254801657
$arr[15,15] = 0.1175616254801657
$arr[15,16] = 0.2798521037145556
$arr[15,17] = 2.518668933431
$arr[15,18] = 0.0001044932585662284
$arr[15,19] = 0.0001044932585662284

$ws.Range("A2:T17").Value2 = $arr